$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Icam1"
$ws.Cells.Item(2,3).Value = "Il2ra"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 99.11651100000002
$ws.Cells.Item(2,8).Value = 297.3495330000001
$ws.Cells.Item(2,9).Value = 0.799346251215574
$ws.Cells.Item(2,10).Value = 0.7993462512155741
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.3213843333333333
$ws.Cells.Item(2,14).Value = 0.964153
$ws.Cells.Item(2,15).Value = 0.1066545648432073
$ws.Cells.Item(2,16).Value = 0.1066545648432073
$ws.Cells.Item(2,17).Value = 31.854493810061
$ws.Cells.Item(2,18).Value = 286.6904442905491
$ws.Cells.Item(2,19).Value = 0.08525392658244609
$ws.Cells.Item(2,20).Value = 0.08525392658244609

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Icam1"
$ws.Cells.Item(3,3).Value = "Il2ra"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 99.11651100000002
$ws.Cells.Item(3,8).Value = 297.3495330000001
$ws.Cells.Item(3,9).Value = 0.799346251215574
$ws.Cells.Item(3,10).Value = 0.7993462512155741
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.174648
$ws.Cells.Item(3,14).Value = 3.523944
$ws.Cells.Item(3,15).Value = 0.3898185390200842
$ws.Cells.Item(3,16).Value = 0.3898185390200842
$ws.Cells.Item(3,17).Value = 116.427011413128
$ws.Cells.Item(3,18).Value = 1047.843102718152
$ws.Cells.Item(3,19).Value = 0.3115999878200362
$ws.Cells.Item(3,20).Value = 0.3115999878200363

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Icam1"
$ws.Cells.Item(4,3).Value = "Il2ra"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 99.11651100000002
$ws.Cells.Item(4,8).Value = 297.3495330000001
$ws.Cells.Item(4,9).Value = 0.799346251215574
$ws.Cells.Item(4,10).Value = 0.7993462512155741
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.517287666666667
$ws.Cells.Item(4,14).Value = 4.551863
$ws.Cells.Item(4,15).Value = 0.5035268961367085
$ws.Cells.Item(4,16).Value = 0.5035268961367085
$ws.Cells.Item(4,17).Value = 150.388259703331
$ws.Cells.Item(4,18).Value = 1353.494337329979
$ws.Cells.Item(4,19).Value = 0.4024923368130917
$ws.Cells.Item(4,20).Value = 0.4024923368130917

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Icam1"
$ws.Cells.Item(5,3).Value = "Il2ra"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 23.582852
$ws.Cells.Item(5,8).Value = 70.748556
$ws.Cells.Item(5,9).Value = 0.1901889417714845
$ws.Cells.Item(5,10).Value = 0.1901889417714845
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.3213843333333333
$ws.Cells.Item(5,14).Value = 0.964153
$ws.Cells.Item(5,15).Value = 0.1066545648432073
$ws.Cells.Item(5,16).Value = 0.1066545648432073
$ws.Cells.Item(5,17).Value = 7.579159168118666
$ws.Cells.Item(5,18).Value = 68.212432513068
$ws.Cells.Item(5,19).Value = 0.02028451882262776
$ws.Cells.Item(5,20).Value = 0.02028451882262776

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Icam1"
$ws.Cells.Item(6,3).Value = "Il2ra"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 23.582852
$ws.Cells.Item(6,8).Value = 70.748556
$ws.Cells.Item(6,9).Value = 0.1901889417714845
$ws.Cells.Item(6,10).Value = 0.1901889417714845
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.174648
$ws.Cells.Item(6,14).Value = 3.523944
$ws.Cells.Item(6,15).Value = 0.3898185390200842
$ws.Cells.Item(6,16).Value = 0.3898185390200842
$ws.Cells.Item(6,17).Value = 27.701549936096
$ws.Cells.Item(6,18).Value = 249.313949424864
$ws.Cells.Item(6,19).Value = 0.07413917541913594
$ws.Cells.Item(6,20).Value = 0.07413917541913594

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Icam1"
$ws.Cells.Item(7,3).Value = "Il2ra"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 23.582852
$ws.Cells.Item(7,8).Value = 70.748556
$ws.Cells.Item(7,9).Value = 0.1901889417714845
$ws.Cells.Item(7,10).Value = 0.1901889417714845
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.517287666666667
$ws.Cells.Item(7,14).Value = 4.551863
$ws.Cells.Item(7,15).Value = 0.5035268961367085
$ws.Cells.Item(7,16).Value = 0.5035268961367085
$ws.Cells.Item(7,17).Value = 35.78197048442533
$ws.Cells.Item(7,18).Value = 322.037734359828
$ws.Cells.Item(7,19).Value = 0.09576524752972077
$ws.Cells.Item(7,20).Value = 0.09576524752972077

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Icam1"
$ws.Cells.Item(8,3).Value = "Il2ra"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.297604333333333
$ws.Cells.Item(8,8).Value = 3.892813
$ws.Cells.Item(8,9).Value = 0.01046480701294141
$ws.Cells.Item(8,10).Value = 0.01046480701294141
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.3213843333333333
$ws.Cells.Item(8,14).Value = 0.964153
$ws.Cells.Item(8,15).Value = 0.1066545648432073
$ws.Cells.Item(8,16).Value = 0.1066545648432073
$ws.Cells.Item(8,17).Value = 0.4170297035987778
$ws.Cells.Item(8,18).Value = 3.753267332389
$ws.Cells.Item(8,19).Value = 0.001116119438133409
$ws.Cells.Item(8,20).Value = 0.001116119438133409

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Icam1"
$ws.Cells.Item(9,3).Value = "Il2ra"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.297604333333333
$ws.Cells.Item(9,8).Value = 3.892813
$ws.Cells.Item(9,9).Value = 0.01046480701294141
$ws.Cells.Item(9,10).Value = 0.01046480701294141
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.174648
$ws.Cells.Item(9,14).Value = 3.523944
$ws.Cells.Item(9,15).Value = 0.3898185390200842
$ws.Cells.Item(9,16).Value = 0.3898185390200842
$ws.Cells.Item(9,17).Value = 1.524228334941333
$ws.Cells.Item(9,18).Value = 13.718055014472
$ws.Cells.Item(9,19).Value = 0.00407937578091195
$ws.Cells.Item(9,20).Value = 0.00407937578091195

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Icam1"
$ws.Cells.Item(10,3).Value = "Il2ra"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.297604333333333
$ws.Cells.Item(10,8).Value = 3.892813
$ws.Cells.Item(10,9).Value = 0.01046480701294141
$ws.Cells.Item(10,10).Value = 0.01046480701294141
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.517287666666667
$ws.Cells.Item(10,14).Value = 4.551863
$ws.Cells.Item(10,15).Value = 0.5035268961367085
$ws.Cells.Item(10,16).Value = 0.5035268961367085
$ws.Cells.Item(10,17).Value = 1.968839051179889
$ws.Cells.Item(10,18).Value = 17.719551460619
$ws.Cells.Item(10,19).Value = 0.005269311793896047
$ws.Cells.Item(10,20).Value = 0.005269311793896047
